$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a "   " (3-space) indentation prefix to the pseudo-code lines that are
# nested inside "Begin if" / "Begin" blocks, effectively tabulating the code.
$ws.Range("B35").Value = "   playerMaxed(playerID) == true"
$ws.Range("B36").Value = "   then handle maxedPlayerException and retrun false"
$ws.Range("B38").Value = "   randomTrainVal <=2"
$ws.Range("B39").Value = "   then return false"
$ws.Range("B41").Value = "   randomTrainVal = 3 or randomTrainVal = 4"
$ws.Range("B42").Value = "   then generate randomHealthLose(1-4) - 1"
$ws.Range("B43").Value = "   addTrainVal(randomTrainVal int, playerID int, position String)"
$ws.Range("B44").Value = "   health = getPlayerHealth(playerID int - randomHealthLose"
$ws.Range("B45").Value = "   getInjured(heath int, playerId int)"
$ws.Range("B46").Value = "   return true"
$ws.Range("B48").Value = "   then generate randomHealthLose(3-4) "
$ws.Range("B49").Value = "   addTrainVal(randomTrainVal int, playerID int, position String)"
$ws.Range("B50").Value = "   health = getPlayerHealth(playerID int - randomHealthLose"
$ws.Range("B51").Value = "   getInjured(heath int, playerId int)"

# Update the saved sheet view/scroll position and selection to match the
# final state left by the author after editing.
$ws.Application.ActiveWindow.ScrollRow = 44
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B51").Select()
